# Auto-generated edit script applying the commit diff to Mateus_Profits workbook
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across 35 rows
# spanning all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 130
$ws.Range("I11").Value = 130
$ws.Range("K11").Value = 130
$ws.Range("M11").Value = 10

$ws.Range("H12").Value = 469.375
$ws.Range("I12").Value = 507.85715
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 507.85715
$ws.Range("L12").Value = 200
$ws.Range("M12").Value = -337.85715
$ws.Range("N12").Value = -540

$ws.Range("H19").Value = 1870.05
$ws.Range("I19").Value = 1699
$ws.Range("J19").Value = 1943.3572
$ws.Range("K19").Value = 1699
$ws.Range("L19").Value = 1943.3572
$ws.Range("M19").Value = -1524
$ws.Range("N19").Value = -2293.3572

$ws.Range("H40").Value = 5319
$ws.Range("I40").Value = 3940
$ws.Range("J40").Value = 6238.3335
$ws.Range("K40").Value = 3940
$ws.Range("L40").Value = 6238.3335
$ws.Range("M40").Value = -3765
$ws.Range("N40").Value = -6588.3335

$ws.Range("H98").Value = 8154.4
$ws.Range("I98").Value = 8154.4
$ws.Range("K98").Value = 8154.4
$ws.Range("M98").Value = -6656.4

$ws.Range("H122").Value = 8154.4
$ws.Range("I122").Value = 8154.4
$ws.Range("K122").Value = 24463.2
$ws.Range("M122").Value = -22013.2

$ws.Range("H131").Value = 83657.53999999999
$ws.Range("I131").Value = 94958.91
$ws.Range("K131").Value = 284876.73
$ws.Range("M131").Value = -279836.73

$ws.Range("H138").Value = 17546420
$ws.Range("I138").Value = 2113.3333
$ws.Range("J138").Value = 23812242
$ws.Range("K138").Value = 6339.999899999999
$ws.Range("L138").Value = 71436726
$ws.Range("M138").Value = -1199.999899999999
$ws.Range("N138").Value = -71447006

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5377.98
$ws.Range("I32").Value = 5185.396
$ws.Range("K32").Value = 5185.396
$ws.Range("M32").Value = -4898.396

$ws.Range("H61").Value = 16136199
$ws.Range("I61").Value = 18524802
$ws.Range("K61").Value = 18524802
$ws.Range("M61").Value = -18524590

$ws.Range("H122").Value = 2190.1333
$ws.Range("I122").Value = 2143.3076
$ws.Range("K122").Value = 6429.9228
$ws.Range("M122").Value = -3979.9228

$ws.Range("H132").Value = 7435.643
$ws.Range("I132").Value = 6835.36
$ws.Range("K132").Value = 20506.08
$ws.Range("M132").Value = -17976.08

$ws.Range("H136").Value = 16136199
$ws.Range("I136").Value = 18524802
$ws.Range("K136").Value = 55574406
$ws.Range("M136").Value = -55571856

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3033.4614
$ws.Range("I134").Value = 3197.4783
$ws.Range("K134").Value = 9592.4349
$ws.Range("M134").Value = -7057.4349

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()  # cell removed in diff (was -96834.5)

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 7292.0625
$ws.Range("I134").Value = 6090.6924
$ws.Range("K134").Value = 18272.0772
$ws.Range("M134").Value = -15737.0772

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 118988.664
$ws.Range("J37").Value = 118988.664
$ws.Range("L37").Value = 356965.992
$ws.Range("N37").Value = -357189.992

$ws.Range("H132").Value = 20001436
$ws.Range("I132").Value = 33334442
$ws.Range("K132").Value = 300009978
$ws.Range("M132").Value = -300007448

$ws.Range("H140").Value = 2049.64
$ws.Range("I140").Value = 1989.2084
$ws.Range("K140").Value = 5967.6252
$ws.Range("M140").Value = -787.6252000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 1062.5454
$ws.Range("I14").Value = 788.8
$ws.Range("J14").Value = 3800
$ws.Range("K14").Value = 788.8
$ws.Range("L14").Value = 3800
$ws.Range("M14").Value = -620.8
$ws.Range("N14").Value = -4136

$ws.Range("H49").Value = 33333
$ws.Range("J49").Value = 33333
$ws.Range("L49").Value = 33333
$ws.Range("N49").Value = -33701

$ws.Range("H126").Value = 3863.739
$ws.Range("I126").Value = 3243.6875
$ws.Range("K126").Value = 9731.0625
$ws.Range("M126").Value = -7261.0625

$ws.Range("H132").Value = 6132.7295
$ws.Range("I132").Value = 4660.231
$ws.Range("K132").Value = 13980.693
$ws.Range("M132").Value = -11450.693

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3777.7778
$ws.Range("I22").Value = 3166.6667
$ws.Range("J22").Value = 5000
$ws.Range("K22").Value = 3166.6667
$ws.Range("L22").Value = 5000
$ws.Range("M22").Value = -2871.6667
$ws.Range("N22").Value = -5590

$ws.Range("H27").Value = 3777.7778
$ws.Range("I27").Value = 3166.6667
$ws.Range("J27").Value = 5000
$ws.Range("K27").Value = 3166.6667
$ws.Range("L27").Value = 5000
$ws.Range("M27").Value = -3059.6667
$ws.Range("N27").Value = -5214

$ws.Range("H61").Value = 112849.445
$ws.Range("I61").Value = 112849.445
$ws.Range("K61").Value = 112849.445
$ws.Range("M61").Value = -112647.445

$ws.Range("H98").Value = 20000
$ws.Range("J98").Value = 20000
$ws.Range("L98").Value = 20000
$ws.Range("N98").Value = -25990  # new cell added in diff

$ws.Range("H113").Value = 112849.445
$ws.Range("I113").Value = 112849.445
$ws.Range("K113").Value = 112849.445
$ws.Range("M113").Value = -110679.445

$ws.Range("H132").Value = 7932.517
$ws.Range("I132").Value = 8186.913
$ws.Range("K132").Value = 24560.739
$ws.Range("M132").Value = -22030.739

$ws.Range("H136").Value = 3341.64
$ws.Range("J136").Value = 6312.5
$ws.Range("L136").Value = 18937.5
$ws.Range("N136").Value = -24037.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 177779060
$ws.Range("I6").Value = 222222580
$ws.Range("K6").Value = 222222580
$ws.Range("M6").Value = -222222465

$ws.Range("H70").Value = 39791.25
$ws.Range("J70").Value = 40475.715
$ws.Range("L70").Value = 40475.715
$ws.Range("N70").Value = -41105.715

$ws.Range("H73").Value = 39791.25
$ws.Range("J73").Value = 40475.715
$ws.Range("L73").Value = 40475.715
$ws.Range("N73").Value = -42659.715

$ws.Range("H113").Value = 599.88
$ws.Range("J113").Value = 717
$ws.Range("L113").Value = 2151
$ws.Range("N113").Value = -6491

$ws.Range("H132").Value = 2495.84
$ws.Range("I132").Value = 2454.3865
$ws.Range("J132").Value = 2799.8333
$ws.Range("K132").Value = 7363.1595
$ws.Range("L132").Value = 8399.499899999999
$ws.Range("M132").Value = -4833.1595
$ws.Range("N132").Value = -13459.4999
